# Updates cryptos list data (prices & volume changes), mirroring the
# upstream GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.979.71"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.973.61"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'354.02"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'112.27"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "'0.562"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.632"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'39.72"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").Value = "'0.0897"
$ws.Range("E11").Value = "  +4.77%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "'19.99"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "'7.97"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "3.449.33"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").Value = "3.012.17"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "'0.994"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "52.113.46"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'14.46"
$ws.Range("E20").Value = "  +6.23%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "'3.31"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").Value = "0.0₃0990"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("D23").Value = "'71.38"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "'270.63"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'2.80"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'0.180"
$ws.Range("E26").Value = "  +9.06%  "
$ws.Range("D27").Value = "'27.73"
$ws.Range("E27").Value = "  +4.24%  "
$ws.Range("D28").Value = "'7.70"
$ws.Range("E28").Value = "  +21.08%  "
$ws.Range("D29").Value = "'0.114"
$ws.Range("E29").Value = "  +27.42%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'10.73"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'37.64"
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("E33").Value = "  +11.54%  "
$ws.Range("D34").Value = "'52.96"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("D35").Value = "'0.0450"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "'1.99"
$ws.Range("E36").Value = "  +4.48%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("D39").Value = "'19.02"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  +3.75%  "
$ws.Range("D42").Value = "'23.94"
$ws.Range("E42").Value = "  +6.14%  "
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D45").Value = "'3.55"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "2.176.17"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'113.96"
$ws.Range("E48").Value = "  -7.34%  "
$ws.Range("D49").Value = "'0.244"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'0.0342"
$ws.Range("E50").Value = "  +5.75%  "
$ws.Range("D51").Value = "'0.942"
$ws.Range("E51").Value = "  -2.54%  "
